$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "日期" (date) column C. This shifts the old "问题"
#    column (D) into C and the old "解答" column (E) into D, matching the
#    new 4-column layout (编号/来源/问题/解答).
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Delete()

# Widen the (new) column C - best effort match for the authored width
# (70.2727272727273 chars under the authoring font's metrics).
$ws.Columns.Item(3).ColumnWidth = 69.57142857142857

# ---------------------------------------------------------------------------
# 2. Fix up existing rows 2-4 (source/question/answer text + wording tweaks)
# ---------------------------------------------------------------------------

# Row 2 - afterpay installments question
$ws.Range("B2").Value = " TIME HONOUR"
$ws.Range("C2").Value = "afterpay支持几期分期？"
$ws.Range("C2").HorizontalAlignment = -4131
$ws.Range("D2").Value = "您好，对于AfterPay的分期，我们只负责功能支持，至于具体的分期策略完全是由AfterPay决定的。通常来说Afterpay支持4期免息，这是也是最常用的分期策略。"

# Row 3 - GooglePay/ApplePay region question
$ws.Range("B3").Value = " TIME HONOUR"
$ws.Range("C3").Value = "Google pay和Apple pay是否也支持加拿大，英国，欧洲地区？"
$ws.Range("C3").HorizontalAlignment = 5
$ws.Range("D3").Value = "您好，我们支持的地区是和GooglePay以及ApplePay官方保持一致的，通常来说除了少部分受制裁地区，大部分国家和地区都支持。具体可参看以下官方文档https://support.google.com/pay/answer/9023773?hl=zh-Hans#zippy=%2C%E5%9C%A8%E7%BA%BF%E4%BB%98%E6%AC%BE%E6%88%96%E5%9C%A8%E5%BA%94%E7%94%A8%E5%86%85%E4%BB%98%E6%AC%BE"

# Row 4 - mc domain trailing slash question
$ws.Range("B4").Value = "香港深辉扬SENHUIYANG"
$ws.Range("C4").Value = "mc后台配置中心，域名管理中的域名末尾带斜杠有影响吗"
$ws.Range("C4").HorizontalAlignment = 5
$ws.Range("D4").Value = "您好，我们不建议带斜杠，及时带上斜杠并不会直接导致错误。因为此处的域名需要保持全局统一，如果此处带了斜杠，其他所有使用到这个域名的地方都需要带上斜杠。"

# ---------------------------------------------------------------------------
# 3. Append the brand-new FAQ rows (5-11)
# ---------------------------------------------------------------------------

# Row 5 - BAIERTE TRADING / dispute fee question
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "BAIERTE TRADING"
$ws.Range("C5").Value = "当发生争议时，如果客户自己取消了拒付，商家还会被收取争议处理费吗？"
$ws.Range("C5").HorizontalAlignment = 5
$ws.Range("D5").Value = "您好，银行对于拒付，只要发生都会收拒付处理费。客户撤销拒付只会影响拒付率的计算，而不影响银行对拒付手续费的收取。"

# Row 6 - BAIERTE TRADING / low success rate question
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "BAIERTE TRADING"
$ws.Range("C6").Value = "afterpay的成功率很低是怎么回事？"
$ws.Range("C6").HorizontalAlignment = 5
$ws.Range("D6").Value = "您好，订单量对于成功率波动有影响，当订单数据量越多，成功率越趋于平衡，此时的订单成功率才具备较高的参考价值。当订单数据量较少时，建议先积累数据量，再观察数据表现。"

# Row 7 - TIME HONOUR / channel blocked question (fill + wrap)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "TIME HONOUR"
$ws.Range("C7").Value = "系统中没有一笔成功交易，是否可以确定是使用的通道不通畅的原因导致的呢？"
$ws.Range("C7").HorizontalAlignment = 5
$ws.Range("C7").WrapText = $true
$ws.Range("D7").Value = "您好，如果通道不通畅的话，交易是无法送入的。如果系统中存在失败交易、过期交易、处理中等状态的交易，说明通道是畅通的。此时无成功交易需要考虑其他潜在的原因。"

# Row 8 - HONGXINYI / test card risk control question
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "HONGXINYI "
$ws.Range("C8").Value = "测试卡测试时 当天支付几次之后就无法再付款，间隔一周也不行，只能加白名单。是否可以取消掉这个风控，因为可能后续用户会有多次付费的情况。  "
$ws.Range("C8").HorizontalAlignment = 5
$ws.Range("D8").Value = "您好，测试卡可以通过加入白名单的形式来避免风控。实际在用户支付过程中，UseePay风控系统认为短期内多次相同金额的重复支付通常是不正常的，会触发拦截，触发条件是：24小时内，同一张卡，同一个邮箱，同一个ip，同金额订单，最多支付三次。此模式足够正常情况下的用户消费，超出3次将触发拦截。"

# Row 9 - HONGXINYI / refund question (default formatting, no alignment override)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "HONGXINYI "
$ws.Range("C9").Value = "退款是原路返回吗？这边后台能看到退款成功，但是实际卡没有收到退款。"
$ws.Range("D9").Value = "您好，退款是3-7天内原路返回，实际处理一般很快，但并非实时到账。没有收到退款的原因通常是银行处理延迟导致。并且退款就有相对应的ARN编码生成，如消费者存在疑问可以提供ARN码以便和银行核实。 "

# Row 10 - LIGHTSPARK / deposit ratio question
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "LIGHTSPARK"
$ws.Range("C10").Value = "为什么账户的保证金和可提现金额比例完全对不上？"
$ws.Range("C10").HorizontalAlignment = 5
$ws.Range("D10").Value = "您好，保证金是按照每笔交易订单的10%收取，已生成结算的订单才会开始收取保证金。所以保证金最终是占比已结算金额的10%，而不是可提现金额的10%，所以二者会有比例差异。"

# Row 11 - SENHUIYANG / settlement timing question - answer goes in column E
#          ("尚未回复" = not yet answered), column D intentionally left blank.
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "SENHUIYANG"
$ws.Range("C11").Value = "支付通道是实时到账还是延时到账？"
$ws.Range("E11").Value = "尚未回复"

# ---------------------------------------------------------------------------
# 4. Sheet view tweaks: drop the frozen/ scrolled "topLeftCell", move the
#    active selection to D12.
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D12").Select()
